$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'66.952.75"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -1.29%  "
$c = $ws.Range("D3")
$c.Value = "'2.605.29"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -0.63%  "
$ws.Range("E4").Value = "  +0.17%  "
$c = $ws.Range("D5")
$c.Value = "'591.74"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.13%  "
$c = $ws.Range("D6")
$c.Value = "'151.67"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -3.04%  "
$ws.Range("E7").Value = "  +0.09%  "
$c = $ws.Range("D8")
$c.Value = "'0.554"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +1.83%  "
$c = $ws.Range("D9")
$c.Value = "'2.604.82"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.76%  "
$c = $ws.Range("D10")
$c.Value = "'0.122"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -2.38%  "
$ws.Range("E11").Value = "  +0.26%  "
$c = $ws.Range("D12")
$c.Value = "'5.13"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -1.98%  "
$c = $ws.Range("D13")
$c.Value = "'0.343"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -4.17%  "
$c = $ws.Range("D14")
$c.Value = "'27.29"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.65%  "
$c = $ws.Range("D15")
$c.Value = "'3.084.20"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -0.43%  "
$c = $ws.Range("D16")
$c.Value = "'0.0000178"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -4.54%  "
$c = $ws.Range("D17")
$c.Value = "'66.925.76"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.57%  "
$c = $ws.Range("D18")
$c.Value = "'2.607.09"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.33%  "
$c = $ws.Range("D19")
$c.Value = "'363.66"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.60%  "
$c = $ws.Range("D20")
$c.Value = "'10.94"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -4.12%  "
$c = $ws.Range("D21")
$c.Value = "'7.32"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -6.07%  "
$ws.Range("E22").Value = "  -0.78%  "
$c = $ws.Range("D23")
$c.Value = "'2.03"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -1.69%  "
$ws.Range("E24").Value = "  +0.02%  "
$c = $ws.Range("D25")
$c.Value = "'9.88"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -1.85%  "
$c = $ws.Range("D26")
$c.Value = "'67.34"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -5.42%  "
$c = $ws.Range("D27")
$c.Value = "'2.737.42"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.45%  "
$c = $ws.Range("D28")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.39%  "
$c = $ws.Range("D29")
$c.Value = "'574.47"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -3.49%  "
$c = $ws.Range("D30")
$c.Value = "'0.0₃0997"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -5.38%  "
$c = $ws.Range("D31")
$c.Value = "'1.37"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -6.01%  "
$c = $ws.Range("D32")
$c.Value = "'7.66"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -4.17%  "
$ws.Range("E33").Value = "  -2.81%  "
$c = $ws.Range("D34")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("E35").Value = "  -9.17%  "
$c = $ws.Range("D36")
$c.Value = "'1.49"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -3.89%  "
$c = $ws.Range("D37")
$c.Value = "'4.83"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -2.80%  "
$c = $ws.Range("D38")
$c.Value = "'156.30"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +1.26%  "
$c = $ws.Range("D39")
$c.Value = "'18.93"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -3.59%  "
$c = $ws.Range("D40")
$c.Value = "'0.364"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -2.20%  "
$c = $ws.Range("D41")
$c.Value = "'5.20"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -4.44%  "
$c = $ws.Range("D42")
$c.Value = "'1.79"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -3.62%  "
$c = $ws.Range("D43")
$c.Value = "'2.52"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -3.92%  "
$c = $ws.Range("D44")
$c.Value = "'40.98"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -1.16%  "
$ws.Range("E45").Value = "  +0.00%  "
$c = $ws.Range("D46")
$c.Value = "'16.38"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.71%  "
$c = $ws.Range("D47")
$c.Value = "'154.91"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -2.33%  "
$c = $ws.Range("D48")
$c.Value = "'0.0₆0285"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -4.09%  "
$c = $ws.Range("D49")
$c.Value = "'3.70"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -1.27%  "
$c = $ws.Range("D50")
$c.Value = "'21.55"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +2.45%  "
$c = $ws.Range("D51")
$c.Value = "'0.615"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -3.03%  "
